# Applies the 24-11-2023 20:45 update to the Algeria Ligue-1 2023-2024 sheet:
#  - swap the two matches that were recorded out of order around 2023-10-06/07
#    (rows 26 & 27, columns F:V only; A:E identifiers stay put)
#  - swap the two matches that were recorded out of order around 2023-11-10/11
#    (rows 35 & 36, columns F:V only; A:E identifiers stay put)
#  - append two newly scraped matches as rows 49 & 50
#
# Notes on this runtime's quirks (discovered empirically):
#  - The Range/Cells ".Value" getter misbehaves when its result is stored in
#    a variable and reused (it can return stale/bogus "Variant Value..."
#    text instead of the real data) - so all reads use ".Value2" instead,
#    which works reliably.
#  - Calling a function with named parameters (e.g. "-RowA 26 -RowB 27")
#    fails to bind the arguments (they come through empty) - so helper
#    functions below are always called with positional arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRow {
    param(
        [int]$RowA,
        [int]$RowB
    )

    for ($col = 6; $col -le 22; $col++) {
        $valA = $ws.Cells.Item($RowA, $col).Value2
        $valB = $ws.Cells.Item($RowB, $col).Value2
        $ws.Cells.Item($RowA, $col).Value2 = $valB
        $ws.Cells.Item($RowB, $col).Value2 = $valA
    }
}

# --- Swap rows 26 / 27 (CR Belouizdad-Khenchela vs El Bayadh-ASO Chlef) ---
Swap-MatchRow 26 27

# --- Swap rows 35 / 36 (Magra-Ben Aknoun vs Khenchela-Biskra) ---
Swap-MatchRow 35 36

# --- Append new row 49: Magra 1-1 El Bayadh ---
$ws.Range("A48:V48").Copy()
$ws.Range("A49:V49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(49, 1).Value2 = 48
$ws.Cells.Item(49, 2).Value2 = "algeria"
$ws.Cells.Item(49, 3).Value2 = "ligue-1"
$ws.Cells.Item(49, 4).Value2 = "2023-2024"
$ws.Cells.Item(49, 5).Value2 = 45254.63541666666
$ws.Cells.Item(49, 6).Value2 = "Magra"
$ws.Cells.Item(49, 7).Value2 = 1
$ws.Cells.Item(49, 8).Value2 = "El Bayadh"
$ws.Cells.Item(49, 9).Value2 = 1
$ws.Cells.Item(49, 10).Value2 = 2.07
$ws.Cells.Item(49, 11).Value2 = "23/11/2023 07:47"
$ws.Cells.Item(49, 12).Value2 = 2.15
$ws.Cells.Item(49, 13).Value2 = "24/11/2023 14:19"
$ws.Cells.Item(49, 14).Value2 = 2.88
$ws.Cells.Item(49, 15).Value2 = "23/11/2023 07:47"
$ws.Cells.Item(49, 16).Value2 = 2.9
$ws.Cells.Item(49, 17).Value2 = "24/11/2023 15:12"
$ws.Cells.Item(49, 18).Value2 = 4.25
$ws.Cells.Item(49, 19).Value2 = "23/11/2023 07:47"
$ws.Cells.Item(49, 20).Value2 = 4.06
$ws.Cells.Item(49, 21).Value2 = "24/11/2023 14:19"
$ws.Cells.Item(49, 22).Value2 = "https://www.betexplorer.com/football/algeria/ligue-1/magra-el-bayadh/jog1n073/"

# --- Append new row 50: Paradou 0-0 Constantine ---
$ws.Range("A49:V49").Copy()
$ws.Range("A50:V50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(50, 1).Value2 = 49
$ws.Cells.Item(50, 2).Value2 = "algeria"
$ws.Cells.Item(50, 3).Value2 = "ligue-1"
$ws.Cells.Item(50, 4).Value2 = "2023-2024"
$ws.Cells.Item(50, 5).Value2 = 45254.63541666666
$ws.Cells.Item(50, 6).Value2 = "Paradou"
$ws.Cells.Item(50, 7).Value2 = 0
$ws.Cells.Item(50, 8).Value2 = "Constantine"
$ws.Cells.Item(50, 9).Value2 = 0
$ws.Cells.Item(50, 10).Value2 = 1.79
$ws.Cells.Item(50, 11).Value2 = "23/11/2023 07:47"
$ws.Cells.Item(50, 12).Value2 = 1.98
$ws.Cells.Item(50, 13).Value2 = "24/11/2023 15:00"
$ws.Cells.Item(50, 14).Value2 = 3.27
$ws.Cells.Item(50, 15).Value2 = "23/11/2023 07:47"
$ws.Cells.Item(50, 16).Value2 = 3.08
$ws.Cells.Item(50, 17).Value2 = "24/11/2023 15:00"
$ws.Cells.Item(50, 18).Value2 = 4.96
$ws.Cells.Item(50, 19).Value2 = "23/11/2023 07:47"
$ws.Cells.Item(50, 20).Value2 = 4.37
$ws.Cells.Item(50, 21).Value2 = "24/11/2023 15:00"
$ws.Cells.Item(50, 22).Value2 = "https://www.betexplorer.com/football/algeria/ligue-1/paradou-constantine/dOgcmthc/"
